$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "582.00", "2.994.13") are preserved exactly as typed, not
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.062.17'
$ws.Range("E2").Value = '  -2.70%  '
$ws.Range("D3").Value = '2.994.13'
$ws.Range("E3").Value = '  -2.69%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '582.00'
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").Value = '145.35'
$ws.Range("E6").Value = '  -6.36%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.524'
$ws.Range("E8").Value = '  -2.79%  '
$ws.Range("D9").Value = '2.989.11'
$ws.Range("E9").Value = '  -2.96%  '
$ws.Range("D10").Value = '0.148'
$ws.Range("E10").Value = '  -5.63%  '
$ws.Range("D11").Value = '5.79'
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").Value = '0.455'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  -4.60%  '
$ws.Range("D14").Value = '34.30'
$ws.Range("E14").Value = '  -7.01%  '
$ws.Range("D15").Value = '0.123'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").Value = '3.492.10'
$ws.Range("E16").Value = '  -2.61%  '
$ws.Range("D17").Value = '7.07'
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").Value = '62.067.31'
$ws.Range("E18").Value = '  -2.65%  '
$ws.Range("D19").Value = '2.998.69'
$ws.Range("E19").Value = '  -2.64%  '
$ws.Range("D20").Value = '459.20'
$ws.Range("E20").Value = '  -5.33%  '
$ws.Range("D21").Value = '13.92'
$ws.Range("E21").Value = '  -4.54%  '
$ws.Range("D22").Value = '0.683'
$ws.Range("E22").Value = '  -3.98%  '
$ws.Range("D23").Value = '7.40'
$ws.Range("E23").Value = '  -2.82%  '
$ws.Range("D24").Value = '81.16'
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("D25").Value = '2.22'
$ws.Range("E25").Value = '  -9.76%  '
$ws.Range("D26").Value = '12.20'
$ws.Range("E26").Value = '  -5.53%  '
$ws.Range("D27").Value = '10.06'
$ws.Range("E27").Value = '  -5.84%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = '2.61'
$ws.Range("E30").Value = '  -3.11%  '
$ws.Range("D31").Value = '7.01'
$ws.Range("E31").Value = '  -6.72%  '
$ws.Range("D32").Value = '2.09'
$ws.Range("E32").Value = '  -7.64%  '
$ws.Range("D33").Value = '28.15'
$ws.Range("E33").Value = '  +3.09%  '
$ws.Range("E34").Value = '  -4.15%  '
$ws.Range("D35").Value = '0.0₃0793'
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("D36").Value = '1.02'
$ws.Range("E36").Value = '  -4.52%  '
$ws.Range("D37").Value = '5.73'
$ws.Range("E37").Value = '  -5.68%  '
$ws.Range("D38").Value = '2.10'
$ws.Range("E38").Value = '  -6.33%  '
$ws.Range("D39").Value = '50.17'
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("D40").Value = '9.14'
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("E41").Value = '  -11.54%  '
$ws.Range("D42").Value = '0.114'
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = '392.24'
$ws.Range("E43").Value = '  -10.89%  '
$ws.Range("D44").Value = '0.0356'
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("D45").Value = '0.272'
$ws.Range("E45").Value = '  -6.41%  '
$ws.Range("D46").Value = '2.725.38'
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("D47").Value = '36.68'
$ws.Range("E47").Value = '  -7.17%  '
$ws.Range("D48").Value = '129.05'
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("D51").Value = '2.18'
$ws.Range("E51").Value = '  -3.26%  '

# Restore the default (unstyled) cell style on column D so we don't
# leave a stray explicit "Text" number format on these cells.
$ws.Range("D2:D51").Style = "Normal"

